# Switch Login URL / Track URL columns from the old Cloudflare Pages
# domain to the new flyqdrone.in custom domain for every order row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDomain = "https://df6fc469.flyq-air.pages.dev"
$newDomain = "https://flyqdrone.in"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($row = 2; $row -le $lastRow; $row++) {
    $loginCell = $ws.Cells.Item($row, 19)  # Column S: Login URL
    $trackCell = $ws.Cells.Item($row, 20)  # Column T: Track URL

    $loginValue = $loginCell.Value2
    $trackValue = $trackCell.Value2

    if ($loginValue -like "$oldDomain*") {
        $loginCell.Value2 = $loginValue -replace [regex]::Escape($oldDomain), $newDomain
    }

    if ($trackValue -like "$oldDomain*") {
        $trackCell.Value2 = $trackValue -replace [regex]::Escape($oldDomain), $newDomain
    }
}
